$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 9).Value = "sd"
$ws.Cells.Item(2, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(17, 9).Value = "qy"
$ws.Cells.Item(17, 10).Value = "Yes-No-Question"
$ws.Cells.Item(46, 9).Value = "sd"
$ws.Cells.Item(46, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(49, 9).Value = "sv"
$ws.Cells.Item(49, 10).Value = "Statement-opinion"
$ws.Cells.Item(55, 9).Value = "aa"
$ws.Cells.Item(55, 10).Value = "Agree/Accept"
$ws.Cells.Item(63, 9).Value = "b"
$ws.Cells.Item(63, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(76, 9).Value = "sd"
$ws.Cells.Item(76, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(85, 9).Value = "aa"
$ws.Cells.Item(85, 10).Value = "Agree/Accept"
$ws.Cells.Item(90, 9).Value = "b"
$ws.Cells.Item(90, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(95, 9).Value = "sv"
$ws.Cells.Item(95, 10).Value = "Statement-opinion"
$ws.Cells.Item(104, 9).Value = "sv"
$ws.Cells.Item(104, 10).Value = "Statement-opinion"
$ws.Cells.Item(105, 9).Value = "aa"
$ws.Cells.Item(105, 10).Value = "Agree/Accept"
$ws.Cells.Item(110, 9).Value = "aa"
$ws.Cells.Item(110, 10).Value = "Agree/Accept"
$ws.Cells.Item(112, 9).Value = "sv"
$ws.Cells.Item(112, 10).Value = "Statement-opinion"
$ws.Cells.Item(114, 9).Value = "b"
$ws.Cells.Item(114, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(119, 9).Value = "sd"
$ws.Cells.Item(119, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(128, 9).Value = "ba"
$ws.Cells.Item(128, 10).Value = "Appreciation"
$ws.Cells.Item(140, 9).Value = "b"
$ws.Cells.Item(140, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(142, 9).Value = "ba"
$ws.Cells.Item(142, 10).Value = "Appreciation"
$ws.Cells.Item(144, 9).Value = "sd"
$ws.Cells.Item(144, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(151, 9).Value = "sd"
$ws.Cells.Item(151, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(162, 9).Value = "ba"
$ws.Cells.Item(162, 10).Value = "Appreciation"
$ws.Cells.Item(175, 9).Value = "aa"
$ws.Cells.Item(175, 10).Value = "Agree/Accept"
$ws.Cells.Item(176, 9).Value = "sd"
$ws.Cells.Item(176, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(181, 9).Value = "ba"
$ws.Cells.Item(181, 10).Value = "Appreciation"
$ws.Cells.Item(186, 9).Value = "sd"
$ws.Cells.Item(186, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(190, 9).Value = "aa"
$ws.Cells.Item(190, 10).Value = "Agree/Accept"
$ws.Cells.Item(205, 9).Value = "sd"
$ws.Cells.Item(205, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(206, 9).Value = "sv"
$ws.Cells.Item(206, 10).Value = "Statement-opinion"
$ws.Cells.Item(210, 9).Value = "sd"
$ws.Cells.Item(210, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(214, 9).Value = "sd"
$ws.Cells.Item(214, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(216, 9).Value = "aa"
$ws.Cells.Item(216, 10).Value = "Agree/Accept"
$ws.Cells.Item(226, 9).Value = "sv"
$ws.Cells.Item(226, 10).Value = "Statement-opinion"
$ws.Cells.Item(227, 9).Value = "sd"
$ws.Cells.Item(227, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(228, 9).Value = "aa"
$ws.Cells.Item(228, 10).Value = "Agree/Accept"
$ws.Cells.Item(229, 9).Value = "ba"
$ws.Cells.Item(229, 10).Value = "Appreciation"
$ws.Cells.Item(231, 9).Value = "sd"
$ws.Cells.Item(231, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(236, 9).Value = "sd"
$ws.Cells.Item(236, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(247, 9).Value = "sd"
$ws.Cells.Item(247, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(248, 9).Value = "aa"
$ws.Cells.Item(248, 10).Value = "Agree/Accept"
$ws.Cells.Item(262, 9).Value = "sd"
$ws.Cells.Item(262, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(266, 9).Value = "sv"
$ws.Cells.Item(266, 10).Value = "Statement-opinion"
$ws.Cells.Item(269, 9).Value = "aa"
$ws.Cells.Item(269, 10).Value = "Agree/Accept"
$ws.Cells.Item(282, 9).Value = "sv"
$ws.Cells.Item(282, 10).Value = "Statement-opinion"
$ws.Cells.Item(284, 9).Value = "ba"
$ws.Cells.Item(284, 10).Value = "Appreciation"
$ws.Cells.Item(287, 9).Value = "b"
$ws.Cells.Item(287, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(299, 9).Value = "ba"
$ws.Cells.Item(299, 10).Value = "Appreciation"
$ws.Cells.Item(318, 9).Value = "%"
$ws.Cells.Item(318, 10).Value = "Uninterpretable"
$ws.Cells.Item(319, 9).Value = "%"
$ws.Cells.Item(319, 10).Value = "Uninterpretable"
$ws.Cells.Item(324, 9).Value = "sd"
$ws.Cells.Item(324, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(350, 9).Value = "%"
$ws.Cells.Item(350, 10).Value = "Uninterpretable"
$ws.Cells.Item(351, 9).Value = "sd"
$ws.Cells.Item(351, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(357, 9).Value = "sd"
$ws.Cells.Item(357, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(359, 9).Value = "sd"
$ws.Cells.Item(359, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(364, 9).Value = "aa"
$ws.Cells.Item(364, 10).Value = "Agree/Accept"
$ws.Cells.Item(366, 9).Value = "sv"
$ws.Cells.Item(366, 10).Value = "Statement-opinion"
$ws.Cells.Item(375, 9).Value = "sd"
$ws.Cells.Item(375, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(381, 9).Value = "sv"
$ws.Cells.Item(381, 10).Value = "Statement-opinion"
$ws.Cells.Item(384, 9).Value = "sv"
$ws.Cells.Item(384, 10).Value = "Statement-opinion"
$ws.Cells.Item(386, 9).Value = "sd"
$ws.Cells.Item(386, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(395, 9).Value = "sd"
$ws.Cells.Item(395, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(396, 9).Value = "sd"
$ws.Cells.Item(396, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(403, 9).Value = "sd"
$ws.Cells.Item(403, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(411, 9).Value = "aa"
$ws.Cells.Item(411, 10).Value = "Agree/Accept"
$ws.Cells.Item(414, 9).Value = "aa"
$ws.Cells.Item(414, 10).Value = "Agree/Accept"
$ws.Cells.Item(417, 9).Value = "%"
$ws.Cells.Item(417, 10).Value = "Uninterpretable"
$ws.Cells.Item(418, 9).Value = "sv"
$ws.Cells.Item(418, 10).Value = "Statement-opinion"
$ws.Cells.Item(421, 9).Value = "%"
$ws.Cells.Item(421, 10).Value = "Uninterpretable"
$ws.Cells.Item(425, 9).Value = "sv"
$ws.Cells.Item(425, 10).Value = "Statement-opinion"
$ws.Cells.Item(426, 9).Value = "sd"
$ws.Cells.Item(426, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(427, 9).Value = "sv"
$ws.Cells.Item(427, 10).Value = "Statement-opinion"
$ws.Cells.Item(434, 9).Value = "sd"
$ws.Cells.Item(434, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(462, 9).Value = "sv"
$ws.Cells.Item(462, 10).Value = "Statement-opinion"
$ws.Cells.Item(464, 9).Value = "sv"
$ws.Cells.Item(464, 10).Value = "Statement-opinion"
$ws.Cells.Item(465, 9).Value = "sv"
$ws.Cells.Item(465, 10).Value = "Statement-opinion"
$ws.Cells.Item(483, 9).Value = "sv"
$ws.Cells.Item(483, 10).Value = "Statement-opinion"
$ws.Cells.Item(499, 9).Value = "sv"
$ws.Cells.Item(499, 10).Value = "Statement-opinion"
$ws.Cells.Item(508, 9).Value = "%"
$ws.Cells.Item(508, 10).Value = "Uninterpretable"
$ws.Cells.Item(516, 9).Value = "aa"
$ws.Cells.Item(516, 10).Value = "Agree/Accept"
$ws.Cells.Item(522, 9).Value = "ba"
$ws.Cells.Item(522, 10).Value = "Appreciation"
$ws.Cells.Item(523, 9).Value = "b"
$ws.Cells.Item(523, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(524, 9).Value = "sv"
$ws.Cells.Item(524, 10).Value = "Statement-opinion"
$ws.Cells.Item(531, 9).Value = "ba"
$ws.Cells.Item(531, 10).Value = "Appreciation"
$ws.Cells.Item(534, 9).Value = "sd"
$ws.Cells.Item(534, 10).Value = "Statement-non-opinion"
